$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 4
$ws.Range("B102").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C102").Value = "Los Lagos"
$ws.Range("D102").Value = 44474
$ws.Range("E102").Value = 10
$ws.Range("F102").Value = 100112039
$ws.Range("G102").Value = "Ciboulette"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 280
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 3000
$ws.Range("N102").Value = "$/docena de atados"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("P102").Value = 1000
$ws.Range("Q102").Value = 3
$ws.Range("R102").Value = "Hortaliza"
